# Update the cryptocurrency price/volume table (columns D and E, rows 2-51).
# Note: some new "Price" values are plain numeric-looking strings (e.g. "427.30").
# A leading apostrophe forces Excel to store them as text (matching the
# original inline-string cell type) instead of coercing them to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.431.18'
$ws.Range("E2").Value = '  +4.66%  '
$ws.Range("D3").Value = '3.844.11'
$ws.Range("E3").Value = '  +9.41%  '
$ws.Range("D4").Value = '''1.00'
$ws.Range("E4").Value = '  -0.33%  '
$ws.Range("D5").Value = '''427.30'
$ws.Range("E5").Value = '  +9.21%  '
$ws.Range("D6").Value = '''131.81'
$ws.Range("E6").Value = '  +8.47%  '
$ws.Range("D7").Value = '3.836.36'
$ws.Range("E7").Value = '  +9.46%  '
$ws.Range("D8").Value = '''0.615'
$ws.Range("E8").Value = '  +5.04%  '
$ws.Range("E9").Value = '  -0.10%  '
$ws.Range("D10").Value = '''0.732'
$ws.Range("E10").Value = '  +8.63%  '
$ws.Range("D11").Value = '''0.158'
$ws.Range("E11").Value = '  +4.76%  '
$ws.Range("D12").Value = '''0.0000342'
$ws.Range("E12").Value = '  +0.24%  '
$ws.Range("D13").Value = '''41.97'
$ws.Range("E13").Value = '  +8.49%  '
$ws.Range("E14").Value = '  +13.93%  '
$ws.Range("D15").Value = '4.448.10'
$ws.Range("E15").Value = '  +9.26%  '
$ws.Range("D16").Value = '''15.92'
$ws.Range("E16").Value = '  +26.56%  '
$ws.Range("D17").Value = '3.892.30'
$ws.Range("E17").Value = '  +10.51%  '
$ws.Range("E18").Value = '  +1.34%  '
$ws.Range("D19").Value = '''20.08'
$ws.Range("E19").Value = '  +7.13%  '
$ws.Range("E20").Value = '  +8.48%  '
$ws.Range("D21").Value = '66.616.06'
$ws.Range("E21").Value = '  +4.45%  '
$ws.Range("D22").Value = '''416.14'
$ws.Range("E22").Value = '  +6.06%  '
$ws.Range("D23").Value = '''15.12'
$ws.Range("E23").Value = '  +9.39%  '
$ws.Range("D24").Value = '''85.13'
$ws.Range("E24").Value = '  +5.67%  '
$ws.Range("D25").Value = '''3.11'
$ws.Range("E25").Value = '  +8.65%  '
$ws.Range("D26").Value = '''37.64'
$ws.Range("E26").Value = '  +14.00%  '
$ws.Range("D27").Value = '''10.05'
$ws.Range("E27").Value = '  +14.81%  '
$ws.Range("E28").Value = '  +9.71%  '
$ws.Range("D29").Value = '''5.35'
$ws.Range("E29").Value = '  +2.17%  '
$ws.Range("D30").Value = '''9.27'
$ws.Range("E30").Value = '  +36.00%  '
$ws.Range("D31").Value = '''720.64'
$ws.Range("E31").Value = '  +7.74%  '
$ws.Range("D32").Value = '''13.84'
$ws.Range("E32").Value = '  +15.79%  '
$ws.Range("D33").Value = '''0.126'
$ws.Range("E33").Value = '  +15.13%  '
$ws.Range("D34").Value = '''2.78'
$ws.Range("E34").Value = '  +6.30%  '
$ws.Range("D35").Value = '''1.00'
$ws.Range("E35").Value = '  -0.11%  '
$ws.Range("D36").Value = '''5.80'
$ws.Range("E36").Value = '  +44.48%  '
$ws.Range("D37").Value = '''39.28'
$ws.Range("E37").Value = '  +7.10%  '
$ws.Range("D38").Value = '''0.152'
$ws.Range("E38").Value = '  +1.14%  '
$ws.Range("D39").Value = '''55.75'
$ws.Range("E39").Value = '  +3.27%  '
$ws.Range("D40").Value = '0.0₃0747'
$ws.Range("E40").Value = '  +18.44%  '
$ws.Range("D41").Value = '''0.0466'
$ws.Range("E41").Value = '  +6.90%  '
$ws.Range("D42").Value = '''2.91'
$ws.Range("E42").Value = '  +6.79%  '
$ws.Range("E43").Value = '  +0.23%  '
$ws.Range("D44").Value = '''3.27'
$ws.Range("E44").Value = '  +6.76%  '
$ws.Range("E45").Value = '  +4.82%  '
$ws.Range("E46").Value = '  +10.64%  '
$ws.Range("D47").Value = '''0.321'
$ws.Range("E47").Value = '  +15.78%  '
$ws.Range("D48").Value = '''2.87'
$ws.Range("E48").Value = '  +6.03%  '
$ws.Range("D49").Value = '''142.64'
$ws.Range("E49").Value = '  +2.44%  '
$ws.Range("E50").Value = '  +5.02%  '
$ws.Range("D51").Value = '''2.06'
$ws.Range("E51").Value = '  +5.91%  '
